$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-04-29"

# Update the column header label for the April 2022 column (B)
$ws.Range("B1").Value = "April 2022 (through April 29)"

# Update existing counts that changed
$ws.Range("B2").Value = 10    # Austin
$ws.Range("F3").Value = 2     # Englewood
$ws.Range("AD3").Value = 2    # Englewood
$ws.Range("R4").Value = 2     # North Lawndale
$ws.Range("J5").Value = 4     # Garfield Park
$ws.Range("B48").Value = 2    # Grand Boulevard

# New data points added for 2022-05-07
$ws.Range("V13").Value = 1    # Wicker Park
$ws.Range("N14").Value = 1    # Belmont Cragin
$ws.Range("F21").Value = 1    # Auburn Gresham
$ws.Range("B27").Value = 1    # West Loop
$ws.Range("J72").Value = 1    # Hermosa
